$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Currency-style number format used for the "Total Spent" column (C),
# matching numFmtId 8 ("$#,##0.00_);[Red]($#,##0.00)").
$currencyFormat = '"$"#,##0.00_);[Red]\("$"#,##0.00\)'

$rows = @(
    @{ Date = "3/1/2020'";  Restaurant = "Deli Zone";           Spent = 8.47;   Rewards = 1 },
    @{ Date = "3/3/2020'";  Restaurant = "Chipotle";             Spent = 11.75;  Rewards = 2 },
    @{ Date = "3/10/2020'"; Restaurant = "Cheeba Hut";           Spent = 9.97;   Rewards = 1 },
    @{ Date = "3/16/2020'"; Restaurant = "Flagstaff House";      Spent = 178.67; Rewards = 5 },
    @{ Date = "3/18/2020'"; Restaurant = "Moto Maki";            Spent = 12.42;  Rewards = 2 },
    @{ Date = "3/20/2020'"; Restaurant = "Illegal Peet's";       Spent = 10.51;  Rewards = 1 },
    @{ Date = "3/25/2020'"; Restaurant = "The Buff";             Spent = 28.3;   Rewards = 2 },
    @{ Date = "3/30/2020'"; Restaurant = "Brasserie ten ten";    Spent = 65.45;  Rewards = 1 },
    @{ Date = "4/1/2020'";  Restaurant = "Rio";                  Spent = 80.65;  Rewards = 3 },
    @{ Date = "4/2/2020'";  Restaurant = "The West End Tavern";  Spent = 45.62;  Rewards = 1 }
)

$startRow = 5
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Range("A$r").Value = $row.Date
    $ws.Range("B$r").Value = $row.Restaurant

    $cell = $ws.Range("C$r")
    $cell.Value = $row.Spent
    $cell.NumberFormat = $currencyFormat

    $ws.Range("D$r").Value = $row.Rewards
}

$ws.Columns.Item(4).ColumnWidth = 14.166666666666666
$ws.Columns.Item(5).ColumnWidth = 15.166666666666666

$ws.Range("E12").Select() | Out-Null
